$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows that correspond to internship/job postings that were removed
# (delete from the bottom up so row numbers of earlier rows stay valid)
$rowsToDelete = @(25, 23, 21, 19, 17, 14, 8, 6, 5, 4, 3)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# Append the two new rows for "System Technic Engineering Pte Ltd"
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow1 = $lastRow + 1
$newRow2 = $lastRow + 2

$ws.Cells.Item($newRow1, 1).Value = "System Technic Engineering Pte Ltd"
$ws.Cells.Item($newRow1, 2).Value = "Full-Time"
$ws.Cells.Item($newRow1, 3).Value = "Less than a year"
$ws.Cells.Item($newRow1, 4).Value = "https://glints.com/sg/opportunities/jobs/engineering-intern/d1f645f0-696d-4682-bc47-1a174108c33a?utm_referrer=explore"

$ws.Cells.Item($newRow2, 1).Value = "System Technic Engineering Pte Ltd"
$ws.Cells.Item($newRow2, 2).Value = "Full-Time"
$ws.Cells.Item($newRow2, 3).Value = "1 – 3 years"
$ws.Cells.Item($newRow2, 4).Value = "https://glints.com/sg/opportunities/jobs/administration-assistant/515c46cf-61dd-424e-8733-029faa22556d?utm_referrer=explore"

# Re-sort the data range (A2:D<last>) by Company, Job Type, Experience - matching the
# workbook's existing sortState (A1, B1, C1 as the key headers)
$dataRange = $ws.Range("A2:D" + $newRow2)
$key1 = $ws.Range("A1")
$key2 = $ws.Range("B1")
$key3 = $ws.Range("C1")
$dataRange.Sort($key1, 1, $key2, 0, 1, $key3, 1, 0)

# Widen column A to fit the longer company name now present in the data
$ws.Columns.Item(1).ColumnWidth = 29.88671875
